$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1, copying the format of H1 (bold, centered, bordered header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-23
$values = @(
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(7,8),
    @(7,7),
    @(9,9),
    @(7,7),
    @(8,8),
    @(7,7),
    @(9,9),
    @(8,8),
    @(6,6),
    @(6,6),
    @(4,4),
    @(4,4),
    @(5,5),
    @(5,5)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
